# Add a new person row (Vince White) to the People sheet, including a
# hyperlinked "education" cell, matching the "Can replace existing when
# importing people" import behaviour.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row (row 3)
$ws.Range("A3").Value = "Vince White"
$ws.Range("D3").Value = "male"
$ws.Range("E3").Value = 1960

# Add the hyperlink for the education cell first (this is what introduces the
# Hyperlink style/font), then restore the friendly display text in the cell.
$ws.Hyperlinks.Add($ws.Range("F3"), "https://en.wikipedia.org/wiki/University_College_London", "", "University College London", "https://en.wikipedia.org/wiki/University_College_London") | Out-Null
$ws.Range("F3").Value = "University College London"

$ws.Range("H3").Value = "musician"
$ws.Range("I3").Value = "lead guitar"

# Page setup: portrait orientation
$ws.PageSetup.Orientation = 1

# Match the saved selection state
$ws.Range("N10").Select() | Out-Null
